$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: Create Tag Controller (0.25 hrs, Implementation, reuse "basically a copy of above, tested." comment)
$ws.Range("A17").Value = 43504
$ws.Range("A17").NumberFormat = "DD/MM/YY"
$ws.Range("B17").Value = 0.25
$ws.Range("C17").Value = "Implementation"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "Create Tag Controller"
$ws.Range("F17").Value = "basically a copy of above, tested."

# Row 18: Create Transaction Functionality (1 hr, Implementation, new comment about db naming)
$ws.Range("A18").Value = 43504
$ws.Range("A18").NumberFormat = "DD/MM/YY"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "Implementation"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Create Transaction Functionality"
$ws.Range("F18").Value = "Little bit of bother remembering which name points to the db"

$ws.Range("F18").Select()
